$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model/price data for rows 2-11 (B = model, C = price).
# Prices are written as text (NumberFormat "@") so the numeric-looking
# strings are preserved as text rather than being coerced into numbers.
$data = @(
    @{ Row = 2;  Model = "Smartphone Galaxy A04e"; Price = "854.91" },
    @{ Row = 3;  Model = "Smartphone Galaxy A22"; Price = "1349.01" },
    @{ Row = 4;  Model = "Smartphone Multi G Max 2"; Price = "971.91" },
    @{ Row = 5;  Model = "Smartphone Multi G Max 2"; Price = "809.91" },
    @{ Row = 6;  Model = "Smartphone Multi G Max 2"; Price = "809.91" },
    @{ Row = 7;  Model = "Smartphone Multi F"; Price = "533.61" },
    @{ Row = 8;  Model = "Smartphone Multi G 2"; Price = "728.91" },
    @{ Row = 9;  Model = "Smartphone Multilaser G 32gb 5mp 5.5PT P9132 Multi CX 1 UN"; Price = "614.61" },
    @{ Row = 10; Model = "Smartphone Nokia G21 Android 11"; Price = "1748.90" },
    @{ Row = 11; Model = "Smartphone Multi F Pro"; Price = "533.61" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.Model

    $priceCell = $ws.Range("C$r")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $entry.Price
}
